# Atualizado por script em 21-11-2023 14:30
#
# 1) Two pairs of existing fixture rows had their match data swapped
#    (rows 61<->62 and rows 63<->64 - everything except the "Indice"
#    column A, which stays tied to the row position).
# 2) Two new fixture rows (124 and 125) were appended at the end of the
#    sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap content of rows 61/62 and 63/64 (columns F:V; B:E and K/O/S
#        are identical between the two rows of each pair anyway) --------
$r61 = $ws.Range("F61:V61")
$r62 = $ws.Range("F62:V62")
$v61 = $r61.Value2
$v62 = $r62.Value2
$r61.Value2 = $v62
$r62.Value2 = $v61

$r63 = $ws.Range("F63:V63")
$r64 = $ws.Range("F64:V64")
$v63 = $r63.Value2
$v64 = $r64.Value2
$r63.Value2 = $v64
$r64.Value2 = $v63

# --- 2) Append two new rows (124 & 125) at the bottom ------------------
# Row 124
$ws.Range("A123:V123").Copy($ws.Range("A124:V124"))
$ws.Range("A124").Value2 = 123
$ws.Range("B124").Value2 = "serbia"
$ws.Range("C124").Value2 = "prva-liga"
$ws.Range("D124").Value2 = "2023-2024"
$ws.Range("E124").Value2 = 45251.54166666666
$ws.Range("F124").Value2 = "Radnicki Beograd"
$ws.Range("G124").Value2 = 2
$ws.Range("H124").Value2 = "Mladost GAT"
$ws.Range("I124").Value2 = 2
$ws.Range("J124").Value2 = 2.68
$ws.Range("K124").Value2 = "26/09/2023 03:12"
$ws.Range("L124").Value2 = 2.82
$ws.Range("M124").Value2 = "21/11/2023 12:57"
$ws.Range("N124").Value2 = 2.79
$ws.Range("O124").Value2 = "26/09/2023 03:12"
$ws.Range("P124").Value2 = 2.68
$ws.Range("Q124").Value2 = "21/11/2023 12:52"
$ws.Range("R124").Value2 = 2.47
$ws.Range("S124").Value2 = "26/09/2023 03:12"
$ws.Range("T124").Value2 = 2.65
$ws.Range("U124").Value2 = "21/11/2023 12:57"
$ws.Range("V124").Value2 = "https://www.betexplorer.com/football/serbia/prva-liga/radnicki-beograd-mladost-gat/xltKELuq/"

# Row 125
$ws.Range("A124:V124").Copy($ws.Range("A125:V125"))
$ws.Range("A125").Value2 = 124
$ws.Range("B125").Value2 = "serbia"
$ws.Range("C125").Value2 = "prva-liga"
$ws.Range("D125").Value2 = "2023-2024"
$ws.Range("E125").Value2 = 45251.625
$ws.Range("F125").Value2 = "Dubocica"
$ws.Range("G125").Value2 = 2
$ws.Range("H125").Value2 = "Smederevo"
$ws.Range("I125").Value2 = 1
$ws.Range("J125").Value2 = 1.98
$ws.Range("K125").Value2 = "26/09/2023 05:12"
$ws.Range("L125").Value2 = 3.09
$ws.Range("M125").Value2 = "21/11/2023 14:59"
$ws.Range("N125").Value2 = 2.91
$ws.Range("O125").Value2 = "26/09/2023 05:12"
$ws.Range("P125").Value2 = 2.83
$ws.Range("Q125").Value2 = "21/11/2023 14:58"
$ws.Range("R125").Value2 = 3.35
$ws.Range("S125").Value2 = "26/09/2023 05:12"
$ws.Range("T125").Value2 = 2.17
$ws.Range("U125").Value2 = "21/11/2023 14:59"
$ws.Range("V125").Value2 = "https://www.betexplorer.com/football/serbia/prva-liga/dubocica-smederevo/jipODuek/"
